$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B to hold the sample-size ("n") values,
# shifting the existing percent columns (B,C,D -> C,D,E).
$ws.Columns("B").Insert()

# Header for the new column.
$ws.Range("B1").Value = "n"

# Sample-size values per guild row.
$ws.Range("B2").Value = 57
$ws.Range("B3").Value = 184
$ws.Range("B4").Value = 23
$ws.Range("B5").Value = 57
$ws.Range("B6").Value = 120
$ws.Range("B7").Value = 38

# Re-apply the percent number format to the shifted percentage columns so
# they collapse back onto the single shared "Percent" style instead of a
# duplicate style entry.
$ws.Range("C2:E7").NumberFormat = "0%"

# Match the saved selection state.
$ws.Range("E12").Select() | Out-Null
